# Rotate the "Recorded By" (column G) comma-separated list of names left by
# one position for every data row: the first name moves to the end of the
# list, e.g. "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $parts = $raw.Split(",")

    if ($parts.Count -lt 2) {
        continue
    }

    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $rotated = @()
    for ($i = 1; $i -lt $trimmed.Count; $i++) {
        $rotated += $trimmed[$i]
    }
    $rotated += $trimmed[0]

    $joined = [string]::Join(", ", $rotated)
    $cell.Value = $joined
}
